$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# 1. Copy the charging-point design-plan table from Output (A12:F17) down into
#    Input (A13:F18) - same values/styles, one row lower on the new sheet.
$wsOutput.Range("A12:F17").Copy($wsInput.Range("A13:F18"))

# 2. Widen column F on both sheets to fit the new content.
$wsInput.Columns.Item(6).ColumnWidth = 18.666666666666668
$wsOutput.Columns.Item(6).ColumnWidth = 20.666666666666668

# 3. Input becomes the active sheet/tab, with the newly added table selected.
$wsInput.Activate()
$wsInput.Range("A13:F18").Select()

# 4. Output is no longer the active sheet; point its own selection at the
#    (still present) design-plan table too.
$wsOutput.Range("A12:F17").Select()
$wsInput.Activate()
